$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Remove the existing hyperlinks (J4:J7 mailto: links) before anything else.
$ws.Hyperlinks.Delete()

# 2. Clear the sample data rows (4-7) but keep their formatting/styles intact.
$ws.Range("B4:M7").ClearContents()

# 3. Give every data row a uniform 30pt custom height, and (by touching row 391)
#    grow the sheet's used range down through row 391 like the target file.
$ws.Range("M391").Value = "touch"
$ws.Range("M391").ClearContents()
$ws.Rows("4:391").RowHeight = 30

# 4. Apply an AutoFilter on the header row, filtered on the "Aff. Type" column
#    (column M, the 13th field in the A3:M3 range) using a discrete value list
#    so it serializes as <filters><filter val="..."/></filters>.
$filterRange = $ws.Range("A3:M3")
$criteria = @("Aff. Type (Permanent, Contract, Visiting)")
$filterRange.AutoFilter(13, $criteria, 7)

# 5. Register the hidden _FilterDatabase defined name that Excel creates
#    whenever a worksheet has an active AutoFilter.
$ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$3:`$M`$3")
$fdb = $ws.Names.Item("_xlnm._FilterDatabase")
$fdb.Visible = $false

# 6. Move the active selection like in the target workbook.
$ws.Range("E7").Select()
